$d = $word.ActiveDocument

# Locate the paragraph that ends the "Aufgabe 3" remarks (the text right
# before the trailing blank/bookmark paragraph) -- new content goes right
# after it. Falling back to "the paragraph before the last one" keeps
# this working even if the wording were ever slightly different.
$anchor = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*Callback Methode*") {
        $anchor = $p
        break
    }
}
if ($anchor -eq $null) {
    $anchor = $d.Paragraphs.Item($d.Paragraphs.Count - 1)
}

# In this runtime, a Range captured on a paragraph keeps pointing at that
# fixed text position even as new paragraphs get inserted after it -- so
# repeatedly calling InsertParagraphAfter() on the SAME anchor range
# always inserts the new paragraph immediately after that anchor,
# pushing every previously inserted paragraph further down the document.
# We exploit that by inserting the four new paragraphs in *reverse* of
# their desired final order.

# A minimal, completely empty <w:p/> (no run at all), added via InsertXML
# -- a plain InsertParagraphAfter() alone leaves a stray empty <w:r/>
# behind, which InsertXML lets us avoid.
$emptyParaXml = "<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'/>"

# 4th (last) new paragraph -- "Die Programme ...". Inserted first so it
# ends up farthest from the anchor, i.e. immediately before the original
# trailing paragraph.
$null = $anchor.Range.InsertParagraphAfter()
$d.Paragraphs.Item($anchor.Index + 1).Range.Text = "Die Programme wurden weitgehend programmiert und die Theorie wurde auch sehr umfangreich beantwortet. "

# 3rd new paragraph -- "Bewertung: ".
$null = $anchor.Range.InsertParagraphAfter()
$d.Paragraphs.Item($anchor.Index + 1).Range.Text = "Bewertung: "

# 2nd new (blank) paragraph.
$null = $anchor.Range.InsertParagraphAfter()
$null = $d.Paragraphs.Item($anchor.Index + 1).Range.InsertXML($emptyParaXml)

# 1st new (blank) paragraph -- ends up immediately after the anchor.
$null = $anchor.Range.InsertParagraphAfter()
$null = $d.Paragraphs.Item($anchor.Index + 1).Range.InsertXML($emptyParaXml)

# The original trailing paragraph (now pushed down to the very end of
# the document) holds only the _GoBack bookmark; prepend the
# "Vorgeschlagenen Note: 8" text to it so the bookmark stays inside the
# same paragraph, right after the new text -- matching the target edit.
$pLast = $d.Paragraphs.Item($d.Paragraphs.Count)
$null = $pLast.Range.InsertBefore("Vorgeschlagenen Note: 8")
